# Auto-generated Excel COM-interop script
# Applies per-cell numeric updates to the Leve profit tables across all 8 sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 45413.293
$ws.Range("J112").Value = 49337.227
$ws.Range("L112").Value = 148011.681
$ws.Range("N112").Value = -150227.681
$ws.Range("H116").Value = 5042.857
$ws.Range("I116").Value = 5042.857
$ws.Range("K116").Value = 5042.857
$ws.Range("M116").Value = -1600.857
$ws.Range("H131").Value = 5833.4287
$ws.Range("I131").Value = 5833.4287
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 17500.2861
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -12460.2861
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 43033.875
$ws.Range("I132").Value = 51090.95
$ws.Range("K132").Value = 153272.85
$ws.Range("M132").Value = -150742.85
$ws.Range("H137").Value = 1299.5625
$ws.Range("I137").Value = 1091.8462
$ws.Range("K137").Value = 3275.5386
$ws.Range("M137").Value = -725.5385999999999
$ws.Range("H138").Value = 4361.877
$ws.Range("J138").Value = 3081.16
$ws.Range("L138").Value = 9243.48
$ws.Range("N138").Value = -19523.48
$ws.Range("H141").Value = 1447.9231
$ws.Range("I141").Value = 1151.9166
$ws.Range("K141").Value = 3455.7498
$ws.Range("M141").Value = 1724.2502

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 159763.73
$ws.Range("I32").Value = 189548.45
$ws.Range("J32").Value = 16255.546
$ws.Range("K32").Value = 189548.45
$ws.Range("L32").Value = 16255.546
$ws.Range("M32").Value = -189261.45
$ws.Range("N32").Value = -16829.546
$ws.Range("H61").Value = 49569.117
$ws.Range("I61").Value = 35674.355
$ws.Range("J61").Value = 85463.914
$ws.Range("K61").Value = 35674.355
$ws.Range("L61").Value = 85463.914
$ws.Range("M61").Value = -35462.355
$ws.Range("N61").Value = -85887.914
$ws.Range("H63").Value = 20953
$ws.Range("I63").Value = 5672.1113
$ws.Range("K63").Value = 5672.1113
$ws.Range("M63").Value = -4986.1113
$ws.Range("H66").Value = 20953
$ws.Range("I66").Value = 5672.1113
$ws.Range("K66").Value = 28360.5565
$ws.Range("M66").Value = -24928.5565
$ws.Range("H122").Value = 1137.7142
$ws.Range("I122").Value = 792.8
$ws.Range("K122").Value = 2378.4
$ws.Range("M122").Value = 71.60000000000036
$ws.Range("H136").Value = 49569.117
$ws.Range("I136").Value = 35674.355
$ws.Range("J136").Value = 85463.914
$ws.Range("K136").Value = 107023.065
$ws.Range("L136").Value = 256391.742
$ws.Range("M136").Value = -104473.065
$ws.Range("N136").Value = -261491.742

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4014
$ws.Range("I86").Value = 1961.2
$ws.Range("J86").Value = 7435.3335
$ws.Range("K86").Value = 1961.2
$ws.Range("L86").Value = 7435.3335
$ws.Range("M86").Value = -838.2
$ws.Range("N86").Value = -9681.333500000001
$ws.Range("H89").Value = 4014
$ws.Range("I89").Value = 1961.2
$ws.Range("J89").Value = 7435.3335
$ws.Range("K89").Value = 9806
$ws.Range("L89").Value = 37176.6675
$ws.Range("M89").Value = -4190
$ws.Range("N89").Value = -48408.6675
$ws.Range("H92").Value = 18825.25
$ws.Range("J92").Value = 18825.25
$ws.Range("L92").Value = 18825.25
$ws.Range("N92").Value = -23817.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 12249.5
$ws.Range("J39").Value = 12249.5
$ws.Range("L39").Value = 12249.5
$ws.Range("N39").Value = -13031.5
$ws.Range("H41").Value = 12000
$ws.Range("I41").Value = 12000
$ws.Range("K41").Value = 12000
$ws.Range("M41").Value = -11572
$ws.Range("H49").Value = 12249.5
$ws.Range("J49").Value = 12249.5
$ws.Range("L49").Value = 12249.5
$ws.Range("N49").Value = -12613.5
$ws.Range("H58").Value = 2037.0571
$ws.Range("I58").Value = 2148.077
$ws.Range("K58").Value = 2148.077
$ws.Range("M58").Value = -1945.077
$ws.Range("H86").Value = 13167.167
$ws.Range("I86").Value = 4651.5
$ws.Range("K86").Value = 4651.5
$ws.Range("M86").Value = -3528.5
$ws.Range("H89").Value = 13167.167
$ws.Range("I89").Value = 4651.5
$ws.Range("K89").Value = 23257.5
$ws.Range("M89").Value = -17641.5
$ws.Range("H132").Value = 45887
$ws.Range("I132").Value = 60991
$ws.Range("K132").Value = 182973
$ws.Range("M132").Value = -180443
$ws.Range("H136").Value = 2037.0571
$ws.Range("I136").Value = 2148.077
$ws.Range("K136").Value = 6444.231000000001
$ws.Range("M136").Value = -3894.231000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11224956
$ws.Range("J4").Value = 1616153.8
$ws.Range("L4").Value = 4848461.4
$ws.Range("N4").Value = -4848685.4
$ws.Range("H15").Value = 570.2
$ws.Range("I15").Value = 367
$ws.Range("J15").Value = 875
$ws.Range("K15").Value = 1101
$ws.Range("L15").Value = 2625
$ws.Range("M15").Value = -961
$ws.Range("N15").Value = -2905
$ws.Range("H21").Value = 398.6316
$ws.Range("I21").Value = 210.23529
$ws.Range("J21").Value = 2000
$ws.Range("K21").Value = 630.70587
$ws.Range("L21").Value = 6000
$ws.Range("M21").Value = -457.70587
$ws.Range("N21").Value = -6346
$ws.Range("H33").Value = 14358274
$ws.Range("I33").Value = 451.66666
$ws.Range("J33").Value = 25126640
$ws.Range("K33").Value = 2709.99996
$ws.Range("L33").Value = 150759840
$ws.Range("M33").Value = -2426.99996
$ws.Range("N33").Value = -150760406
$ws.Range("H49").Value = 4747.5
$ws.Range("J49").Value = 4663.3335
$ws.Range("L49").Value = 13990.0005
$ws.Range("N49").Value = -14302.0005
$ws.Range("H86").Value = 1498.2
$ws.Range("J86").Value = 1496.5
$ws.Range("L86").Value = 4489.5
$ws.Range("N86").Value = -6861.5
$ws.Range("H89").Value = 1498.2
$ws.Range("J89").Value = 1496.5
$ws.Range("L89").Value = 13468.5
$ws.Range("N89").Value = -25324.5
$ws.Range("H96").Value = 7363676.5
$ws.Range("I96").Value = 11771394
$ws.Range("J96").Value = 17479.666
$ws.Range("K96").Value = 35314182
$ws.Range("L96").Value = 52438.99800000001
$ws.Range("M96").Value = -35312123
$ws.Range("N96").Value = -56556.99800000001
$ws.Range("H99").Value = 76936610
$ws.Range("I99").Value = 200003200
$ws.Range("K99").Value = 600009600
$ws.Range("M99").Value = -600007354

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1757.4762
$ws.Range("I132").Value = 1813.2354
$ws.Range("K132").Value = 5439.706200000001
$ws.Range("M132").Value = -2909.706200000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 10000
$ws.Range("I24").Value = 10000
$ws.Range("K24").Value = 10000
$ws.Range("M24").Value = -9770
$ws.Range("H81").Value = 60380.94
$ws.Range("J81").Value = 143339.42
$ws.Range("L81").Value = 286678.84
$ws.Range("N81").Value = -288800.84
$ws.Range("H84").Value = 60380.94
$ws.Range("J84").Value = 143339.42
$ws.Range("L84").Value = 1433394.2
$ws.Range("N84").Value = -1444002.2
$ws.Range("H101").Value = 28749.5
$ws.Range("J101").Value = 28749.5
$ws.Range("L101").Value = 28749.5
$ws.Range("N101").Value = -35239.5
$ws.Range("H113").Value = 342.32
$ws.Range("J113").Value = 355.75
$ws.Range("L113").Value = 1067.25
$ws.Range("N113").Value = -5407.25
$ws.Range("H132").Value = 3015.5
$ws.Range("I132").Value = 3017.8572
$ws.Range("K132").Value = 9053.571599999999
$ws.Range("M132").Value = -6523.571599999999
